$d = $word.ActiveDocument

# 1) Trim the placeholder suffix " (a completar)" from the existing line.
$d.Content.Find.Execute("- Enlace al video explicativo (a completar)", $false, $false, $false, $false, $false, $true, 1, $false, "- Enlace al video explicativo", 2)

# 2) Work on the (still) last paragraph of the document body.
$para = $d.Paragraphs.Last
$r = $para.Range

# Append ": " (colon + space) after the existing text, inheriting the
# surrounding Arial / Courier New / szCs24 run formatting.
$r.InsertAfter(": ")

# Append a one-character placeholder we will turn into the hyperlink; using
# a disposable marker char lets Hyperlinks.Add anchor exactly at the right
# spot without disturbing the text already in the paragraph.
$para = $d.Paragraphs.Last
$r = $para.Range
$r.InsertAfter("#")
$para = $d.Paragraphs.Last
$markerStart = $para.Range.End - 2
$markerRange = $d.Range($markerStart, $markerStart + 1)
$d.Hyperlinks.Add($markerRange, "https://youtu.be/VdU9cLnkqW0", "", "", "https://youtu.be/VdU9cLnkqW0") | Out-Null

# Style + font touch-up on the freshly inserted hyperlink run.
$videoLink = $d.Hyperlinks.Item($d.Hyperlinks.Count)
$videoLinkRange = $videoLink.Range
$videoLinkRange.Style = "Hipervnculo"
$videoLinkRange.Font.Name = "Arial"

# Trailing space after the hyperlink, matching the source formatting again.
$para = $d.Paragraphs.Last
$r = $para.Range
$r.InsertAfter(" ")

# 3) New paragraph for the GitHub repository link.
$para = $d.Paragraphs.Last
$r = $para.Range
$r.InsertParagraphAfter() | Out-Null

$para = $d.Paragraphs.Last
$r = $para.Range
$r.InsertAfter("- Enlace repositorio en GitHub con todo el contenido: ")

$para = $d.Paragraphs.Last
$r = $para.Range
$r.InsertAfter("#")
$para = $d.Paragraphs.Last
$markerStart2 = $para.Range.End - 2
$markerRange2 = $d.Range($markerStart2, $markerStart2 + 1)
$d.Hyperlinks.Add($markerRange2, "https://github.com/CampanaJ/proyect_integrador-AySO", "", "", "https://github.com/CampanaJ/proyect_integrador-AySO") | Out-Null

$repoLink = $d.Hyperlinks.Item($d.Hyperlinks.Count)
$repoLinkRange = $repoLink.Range
$repoLinkRange.Style = "Hipervnculo"
$repoLinkRange.Font.Name = "Arial"

$para = $d.Paragraphs.Last
$r = $para.Range
$r.InsertAfter(" ")
